$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text on Hoja1!A1 ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $ws1.Cells.Item(1,1)
$text = $cellA1.Text
$text = $text.Replace("1000 Bs = 9.76 = 41512.08 pesos", "1000 Bs = 9.85 = 41822.66 pesos")
$text = $text.Replace("41512.08 pesos = 9.71 = 951.62 Bs", "41822.66 pesos = 9.89 = 979.32 Bs")
$cellA1.Value = $text

# --- Update the rate figures on tasas sheet ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 101.5
$ws2.Range("O10").Value = 4245
$ws2.Range("N12").Value = 4230
$ws2.Range("O12").Value = 99.05
